# Mise à jour mapping corps 6e252a7146deaedf3e4787c95fe2bf14b0736235
#
# 1) Bump the "Date" metadata value on the Metadata sheet.
# 2) Remove the "fr-lm-group-de-questionnaires-devaluation.valeur"
#    row from the Elements sheet (row 7), which shifts the
#    "...evaluation" row up from 8 to 7 and shrinks the used range.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2026-02-04T10:58:36+00:00"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("A7:AJ7").EntireRow.Delete()
